$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the session_id value used across all existing data rows (2-15)
#    from 91cddd36-603e-43fd-a265-70bd56024a18 to b064392b-217a-4b17-9893-fc301ca6f256
$newSessionId = "b064392b-217a-4b17-9893-fc301ca6f256"
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newSessionId
}

# 2. Append three new example rows (16-18) describing "Pattern Mismatch" issues
$ws.Cells.Item(16, 1).Value = $newSessionId
$ws.Cells.Item(16, 2).Value = 4
$ws.Cells.Item(16, 3).Value = "Pattern Mismatch"
$ws.Cells.Item(16, 4).Value = "column8"
$ws.Cells.Item(16, 5).Value = "ABC-ABCD"
$ws.Cells.Item(16, 6).Value = "Value ABC-ABCD in column8 does not match the pattern ABC-1234"
$ws.Cells.Item(16, 7).Value = "Follow the pattern ABC-1234 in column8"

$ws.Cells.Item(17, 1).Value = $newSessionId
$ws.Cells.Item(17, 2).Value = 5
$ws.Cells.Item(17, 3).Value = "Pattern Mismatch"
$ws.Cells.Item(17, 4).Value = "column8"
$ws.Cells.Item(17, 5).Value = "XYZ-1234"
$ws.Cells.Item(17, 6).Value = "Value XYZ-1234 in column8 does not match the pattern ABC-1234"
$ws.Cells.Item(17, 7).Value = "Follow the pattern ABC-1234 in column8"

$ws.Cells.Item(18, 1).Value = $newSessionId
$ws.Cells.Item(18, 2).Value = 10
$ws.Cells.Item(18, 3).Value = "Pattern Mismatch"
$ws.Cells.Item(18, 4).Value = "column8"
$ws.Cells.Item(18, 5).Value = "XYZ-1234"
$ws.Cells.Item(18, 6).Value = "Value XYZ-1234 in column8 does not match the pattern ABC-1234"
$ws.Cells.Item(18, 7).Value = "Follow the pattern ABC-1234 in column8"
